# Updated quest items and quests.
# Appends 5 new quest rows (47-51) to the "Quests" sheet, introducing the
# "Creator and the Smith" / "Church curse" quest chain and its items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quests")

# Row 47: The Creator and The Smith (parent quest)
$ws.Range("A47").Value = "The Creator and The Smith"
$ws.Range("B47").Value = "DrunkenAdventurer"
$ws.Range("D47").Value = 15000
$ws.Range("E47").Value = 2000
$ws.Range("I47").Value = 1000
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 1000000
$ws.Range("L47").Value = 25
$ws.Range("O47").Value = 1

# Row 48: Smithies tools of the trade
$ws.Range("A48").Value = "Smithies tools of the trade"
$ws.Range("B48").Value = "Dungeon Master"
$ws.Range("C48").Value = "Smithies Hammer"
$ws.Range("D48").Value = 30000
$ws.Range("E48").Value = 3000
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 2000000
$ws.Range("L48").Value = 25
$ws.Range("P48").Value = "The Creator and The Smith"
$ws.Range("Q48").Value = "Smithies Iron Chunks"

# Row 49: Candles in the dark
$ws.Range("A49").Value = "Candles in the dark"
$ws.Range("B49").Value = "Candle Maker"
$ws.Range("C49").Value = "Broken Smithies Anvil"
$ws.Range("H49").Value = "Candle of the Smithy"
$ws.Range("P49").Value = "Smithies tools of the trade"
$ws.Range("Q49").Value = "Smithies Dying Ember"

# Row 50: The curse of the Church
$ws.Range("A50").Value = "The curse of the Church"
$ws.Range("B50").Value = "Candle Maker"
$ws.Range("C50").Value = "Vial of Water from the well of the Smith"
$ws.Range("D50").Value = 30000
$ws.Range("E50").Value = 5000
$ws.Range("F50").Value = 5000000
$ws.Range("G50").Value = 1000
$ws.Range("H50").Value = "Purgatories Cursed Candle"
$ws.Range("I50").Value = 1500
$ws.Range("J50").Value = 2500
$ws.Range("K50").Value = 2500000
$ws.Range("L50").Value = 100
$ws.Range("P50").Value = "Candles in the dark"
$ws.Range("Q50").Value = "Corrupted Candle of the Church"

# Row 51: Into the House
$ws.Range("A51").Value = "Into the House"
$ws.Range("B51").Value = "Shade Lord"
$ws.Range("C51").Value = "Bag of Transformation"
$ws.Range("D51").Value = 50000
$ws.Range("E51").Value = 5000
$ws.Range("F51").Value = 10000000
$ws.Range("H51").Value = "Purgatory Smiths House Key"
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 1000
$ws.Range("K51").Value = 1000000000
$ws.Range("L51").Value = 100
$ws.Range("P51").Value = "Candles in the dark"

# Column C ("item_id") widened to fit the longest new value
# ("Vial of Water from the well of the Smith").
$ws.Columns.Item(3).ColumnWidth = 47.15
